$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 (pushes "Evolution and Human Behavior" and every
# journal after it down by one row) to make room for the new entry.
$ws.Range("A11:E11").Insert(-4121) | Out-Null

# Match the formatting already used throughout the journal list (row 12, which
# used to be row 11 before the shift).
$ws.Range("A12:E12").Copy() | Out-Null
$ws.Range("A11:E11").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = 0

$ws.Range("E11").Value = "\href{https://royalsocietypublishing.org/journal/rsbl}{Biology Letters}"

$ws.Range("E13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
